# 徐耀昌 2012-04-24 財產申報表 — "汽車" (car) sheet
# Commit: "#5: property boat&car done"
#
# Before this edit the 汽車 (car) sheet was missing its header row (row 1
# held a stray copy of the data instead of field-name labels) and only
# carried 7 of the full 14 columns used by the other property sheets
# (土地/建物). This fills in the proper header labels and extends the
# sheet with the remaining columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) for both the
# header row and the one data record, matching the schema already used
# by the 土地 (land) and 建物 (building) sheets. It also adds a new
# "capacity" shared string for the car's engine displacement column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Extend row 1 / row 2 with the missing columns (H:N) first, by copying
# the existing end-of-row cell's format into each new cell so the new
# cells land in the same visual style family as their row (bold/bordered
# header style for row 1, plain style for row 2) instead of Excel's
# unformatted default.
"H1","I1","J1","K1","L1","M1","N1" | ForEach-Object {
  $ws.Range("G1").Copy($ws.Range($_))
}
"H2","I2","J2","K2","L2","M2","N2" | ForEach-Object {
  $ws.Range("G2").Copy($ws.Range($_))
}

# --- Header row (row 1): field-name labels ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2): the one car record ---
$ws.Range("A2").Value = 40
$ws.Range("B2").Value = "鈴木HMAY1.3LJLXA43D"
$ws.Range("C2").Value = 1328
$ws.Range("D2").Value = "蔡麗卿"
$ws.Range("E2").Value = "94年05月20日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = "(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# J2 ("2012-04-24") must stay plain text (like the rest of the workbook's
# date-ish columns), not get auto-converted into a date serial — force
# Text format for the assignment, then drop back to the row's normal
# style so the cell's effective formatting matches its neighbours.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-24"
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Value = "徐耀昌"
$ws.Range("L2").Value = 921
$ws.Range("M2").Value = "tmp832b1"
$ws.Range("N2").Value = 40
